$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to Text *before* writing so that
# numeric- or percentage-looking strings (e.g. "1.001", "0.000009994",
# "29.190.75", "  -0.55%  ") are stored verbatim as text instead of
# being auto-coerced by Excel into numbers / dates / percentages.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "29.190.75"
$ws.Range("E2").Value = "  -0.55%  "
$ws.Range("D3").Value = "1.825.01"
$ws.Range("E3").Value = "  -0.85%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "235.80"
$ws.Range("E5").Value = "  -1.50%  "
$ws.Range("D6").Value = "0.6101"
$ws.Range("E6").Value = "  -2.98%  "
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("D8").Value = "0.07085"
$ws.Range("E8").Value = "  -4.80%  "
$ws.Range("D9").Value = "0.2802"
$ws.Range("E9").Value = "  -3.12%  "
$ws.Range("D10").Value = "23.46"
$ws.Range("E10").Value = "  -6.03%  "
$ws.Range("D11").Value = "0.07658"
$ws.Range("E11").Value = "  -0.94%  "
$ws.Range("D12").Value = "1.830.24"
$ws.Range("E12").Value = "  -0.61%  "
$ws.Range("D13").Value = "4.805"
$ws.Range("E13").Value = "  -3.18%  "
$ws.Range("D14").Value = "0.000009994"
$ws.Range("E14").Value = "  -2.34%  "
$ws.Range("D15").Value = "0.6313"
$ws.Range("E15").Value = "  -6.47%  "
$ws.Range("D16").Value = "2.066.53"
$ws.Range("E16").Value = "  -1.12%  "
$ws.Range("D17").Value = "78.56"
$ws.Range("E17").Value = "  -3.72%  "
$ws.Range("D18").Value = "5.856"
$ws.Range("E18").Value = "  -6.06%  "
$ws.Range("D19").Value = "29.167.31"
$ws.Range("E19").Value = "  -0.67%  "
$ws.Range("D20").Value = "226.45"
$ws.Range("E20").Value = "  -1.13%  "
$ws.Range("E21").Value = "  +0.14%  "
$ws.Range("D22").Value = "11.76"
$ws.Range("E22").Value = "  -4.24%  "
$ws.Range("D23").Value = "6.999"
$ws.Range("E23").Value = "  -4.65%  "
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("E25").Value = "  -1.61%  "
$ws.Range("D26").Value = "8.041"
$ws.Range("E26").Value = "  -5.08%  "
$ws.Range("D27").Value = "0.1305"
$ws.Range("E27").Value = "  -3.02%  "
$ws.Range("D28").Value = "16.56"
$ws.Range("E28").Value = "  -4.63%  "
$ws.Range("D29").Value = "1.494"
$ws.Range("E29").Value = "  +1.76%  "
$ws.Range("D30").Value = "0.06346"
$ws.Range("E30").Value = "  -13.18%  "
$ws.Range("D31").Value = "1.452"
$ws.Range("E31").Value = "  -1.63%  "
$ws.Range("D32").Value = "3.818"
$ws.Range("E32").Value = "  -5.46%  "
$ws.Range("D33").Value = "3.798"
$ws.Range("E33").Value = "  -5.90%  "
$ws.Range("D34").Value = "1.123"
$ws.Range("E34").Value = "  -1.38%  "
$ws.Range("D35").Value = "1.737"
$ws.Range("E35").Value = "  -4.70%  "
$ws.Range("D36").Value = "0.6452"
$ws.Range("E36").Value = "  -6.97%  "
$ws.Range("D37").Value = "2.547"
$ws.Range("E37").Value = "  -1.23%  "
$ws.Range("D38").Value = "1.214.35"
$ws.Range("E38").Value = "  -1.43%  "
$ws.Range("D39").Value = "2.725"
$ws.Range("E39").Value = "  -2.88%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "6.554"
$ws.Range("E40").Value = "  -4.03%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "0.01737"
$ws.Range("E41").Value = "  -5.53%  "
$ws.Range("D42").Value = "0.9125"
$ws.Range("E42").Value = "  -2.16%  "
$ws.Range("D43").Value = "1.001"
$ws.Range("E43").Value = "  +0.12%  "
$ws.Range("D44").Value = "101.01"
$ws.Range("E44").Value = "  +0.47%  "
$ws.Range("D45").Value = "1.976.31"
$ws.Range("E45").Value = "  -0.74%  "
$ws.Range("D46").Value = "62.60"
$ws.Range("E46").Value = "  -4.16%  "
$ws.Range("E47").Value = "  -4.35%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "8.582"
$ws.Range("E48").Value = "  -3.76%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "1.604"
$ws.Range("E49").Value = "  -6.01%  "
$ws.Range("E50").Value = "  -0.56%  "
$ws.Range("D51").Value = "0.05517"
$ws.Range("E51").Value = "  -2.67%  "

# Restore the default (un-styled) cell style so the written cells keep
# the same "no explicit style index" shape as the rest of the sheet.
$dataRange.Style = "Normal"

